$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.110.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.30%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.228.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.25%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'322.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'98.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -9.33%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.582"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -8.71%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.565"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.57%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -10.44%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'54.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -9.87%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -10.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.566.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -12.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -7.05%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.225.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.93%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.024.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.47%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0968"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -9.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -10.65%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -12.70%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -11.06%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'237.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -10.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.62%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Cosmos"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'10.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -11.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Toncoin"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.25%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -14.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'36.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'20.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -9.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0870"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -9.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'155.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -8.34%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'3.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.62%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'WEMIXToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'2.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -7.30%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -6.94%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -11.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -8.44%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -9.51%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'14.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +8.97%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.07%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.736.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'85.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -13.06%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -11.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -13.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'75.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -12.90%  "
$ws.Range("E51").Style = "Normal"
